$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Update F71: was "2", now "20"
$ws.Cells.Item(71, 6).Value = "'20"

# Add new rows 72-112
$ws.Cells.Item(72, 3).Value = "421_松虫草黑色_scabiosa black_undefined_1bunch"
$ws.Cells.Item(72, 6).Value = "'10"
$ws.Cells.Item(73, 1).Value = "'5"
$ws.Cells.Item(73, 3).Value = "452_粉吊鸟_pink hanging heliconia_undefined_1bunch"
$ws.Cells.Item(73, 6).Value = "'10"
$ws.Cells.Item(74, 1).Value = "'1"
$ws.Cells.Item(74, 3).Value = "522_山归来绿_Smilax china_undefined_1bunch"
$ws.Cells.Item(74, 6).Value = "'5"
$ws.Cells.Item(75, 3).Value = "688_山归来橙_undefined_undefined_1bunch"
$ws.Cells.Item(75, 6).Value = "'5"
$ws.Cells.Item(76, 3).Value = "327_文竹_asparagus fern_undefined_1bunch"
$ws.Cells.Item(76, 6).Value = "'15"
$ws.Cells.Item(77, 3).Value = "532_灯苔_undefined_undefined_1bunch"
$ws.Cells.Item(77, 6).Value = "'20"
$ws.Cells.Item(78, 3).Value = "463_玉兰枝_magnolia flower`nwhite/purple_undefined_1bunch"
$ws.Cells.Item(78, 6).Value = "'10"
$ws.Cells.Item(79, 3).Value = "463_玉兰枝_magnolia flower`nwhite/purple_undefined_1bunch"
$ws.Cells.Item(79, 6).Value = "'10"
$ws.Cells.Item(80, 3).Value = "2_粉洋桔梗_Pink Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Cells.Item(80, 6).Value = "'10"
$ws.Cells.Item(81, 3).Value = "328_卢荀草_undefined_undefined_1bunch"
$ws.Cells.Item(81, 6).Value = "'15"
$ws.Cells.Item(82, 1).Value = "'2"
$ws.Cells.Item(82, 3).Value = "175_火灵鸟_Free Spirit_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(82, 6).Value = "'16"
$ws.Cells.Item(83, 3).Value = "9_茶色洋桔梗_Tea Color Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Cells.Item(83, 6).Value = "'5"
$ws.Cells.Item(84, 3).Value = "614_康乃馨绿_green_undefined_20stems"
$ws.Cells.Item(84, 6).Value = "'15"
$ws.Cells.Item(85, 3).Value = "277_草莓杏仁饼_undefined_Rosa rugosa Thunb._10stems"
$ws.Cells.Item(85, 6).Value = "'15"
$ws.Cells.Item(86, 3).Value = "229_黄蝴蝶_Yellow Butterfly_Rosa rugosa Thunb._10stems"
$ws.Cells.Item(86, 6).Value = "'5"
$ws.Cells.Item(87, 3).Value = "227_多头卡布奇诺_Cappuccino spray_Rosa rugosa Thunb._10stems"
$ws.Cells.Item(87, 6).Value = "'5"
$ws.Cells.Item(88, 3).Value = "550_蕾丝红色_lace flower brown_undefined_1bunch"
$ws.Cells.Item(88, 6).Value = "'10"
$ws.Cells.Item(89, 1).Value = "'3"
$ws.Cells.Item(89, 3).Value = "209_海洋之歌_Ocean Song_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(89, 6).Value = "'5"
$ws.Cells.Item(90, 3).Value = "143_黑巴克_Black Baccara_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(90, 6).Value = "'22"
$ws.Cells.Item(91, 3).Value = "268_猩红泡泡_spray red_Rosa rugosa Thunb._10stems"
$ws.Cells.Item(91, 6).Value = "'5"
$ws.Cells.Item(92, 3).Value = "259_诺拉_Nora_Rosa rugosa Thunb._10stems"
$ws.Cells.Item(92, 6).Value = "'6"
$ws.Cells.Item(93, 3).Value = "238_苏菲宝贝_undefined_Rosa rugosa Thunb._10stems"
$ws.Cells.Item(93, 6).Value = "'10"
$ws.Cells.Item(94, 3).Value = "238_苏菲宝贝_undefined_Rosa rugosa Thunb._10stems"
$ws.Cells.Item(94, 6).Value = "'15"
$ws.Cells.Item(95, 3).Value = "611_康乃馨奶油白_cream white_undefined_20stems"
$ws.Cells.Item(95, 6).Value = "'10"
$ws.Cells.Item(96, 3).Value = "611_康乃馨奶油白_cream white_undefined_20stems"
$ws.Cells.Item(96, 6).Value = "'10"
$ws.Cells.Item(97, 1).Value = "'4"
$ws.Cells.Item(97, 3).Value = "667_大丽花 安吉丽娜_undefined_undefined_5stems"
$ws.Cells.Item(97, 6).Value = "'10"
$ws.Cells.Item(98, 3).Value = "653_大丽花 黑_undefined_undefined_5stems"
$ws.Cells.Item(98, 6).Value = "'5"
$ws.Cells.Item(99, 3).Value = "653_大丽花 黑_undefined_undefined_5stems"
$ws.Cells.Item(99, 6).Value = "'5"
$ws.Cells.Item(100, 3).Value = "677_洋牡丹大香槟_undefined_undefined_1bunch"
$ws.Cells.Item(100, 6).Value = "'10"
$ws.Cells.Item(101, 3).Value = "651_大丽花 奶油桃子_undefined_undefined_5stems"
$ws.Cells.Item(101, 6).Value = "'10"
$ws.Cells.Item(102, 3).Value = "508_风铃花白色_Canterbury Bells `nwhite_undefined_1bunch"
$ws.Cells.Item(102, 6).Value = "'10"
$ws.Cells.Item(103, 1).Value = "'5"
$ws.Cells.Item(103, 3).Value = "414_风铃花粉色_Canterbury Bells`npink_undefined_1bunch"
$ws.Cells.Item(103, 6).Value = "'10"
$ws.Cells.Item(104, 3).Value = "797_维也纳大菊_undefined_undefined_5stems"
$ws.Cells.Item(104, 6).Value = "'20"
$ws.Cells.Item(105, 3).Value = "632_吸色康乃馨紫_tinted purple_undefined_20stems"
$ws.Cells.Item(105, 6).Value = "'10"
$ws.Cells.Item(106, 1).Value = "'6"
$ws.Cells.Item(106, 3).Value = "104_绣球重瓣紫_Hydrangea Purple D_Hydrangea L._1stem"
$ws.Cells.Item(106, 6).Value = "'20"
$ws.Cells.Item(107, 3).Value = "105_绣球莫奈蓝_Hydrangea Monet Blue_Hydrangea L._1stem"
$ws.Cells.Item(107, 6).Value = "'20"
$ws.Cells.Item(108, 3).Value = "148_坦尼克_Tineke_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(108, 6).Value = "'10"
$ws.Cells.Item(109, 3).Value = "83_布拉格_undefined_Gerbera L._10stems"
$ws.Cells.Item(109, 6).Value = "'10"
$ws.Cells.Item(110, 3).Value = "600_康乃馨复古红_vintage red_undefined_20stems"
$ws.Cells.Item(110, 6).Value = "'15"
$ws.Cells.Item(111, 3).Value = "608_康乃馨笑颜_undefined_undefined_20stems"
$ws.Cells.Item(111, 6).Value = "'15"
$ws.Cells.Item(112, 3).Value = "606_康乃馨橙光_orange_undefined_20stems"
$ws.Cells.Item(112, 6).Value = "'10"

# Update Summary sheet G2 value (append additional digits to the existing string)
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Cells.Item(2, 7).Value = "'06101235122013.5105812318129555121010358010105685555851665555538230101010151591310553010155101051010159102020101055152010101015165151555105225610151010105510101010201020201010151510"
